$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H26").Value = 0
$ws.Range("I26").Value = 0
$ws.Range("K26").Value = 0
$ws.Range("M26").ClearContents()

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H38").Value = 1553.7894
$ws.Range("I38").Value = 210.75
$ws.Range("J38").Value = 3856.1428
$ws.Range("K38").Value = 632.25
$ws.Range("L38").Value = 11568.4284
$ws.Range("M38").Value = -260.25
$ws.Range("N38").Value = -12312.4284

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H39").Value = 6175.706
$ws.Range("I39").Value = 311.75
$ws.Range("K39").Value = 935.25
$ws.Range("M39").Value = -639.25

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H51").Value = 3999.75
$ws.Range("J51").Value = 1999.5
$ws.Range("L51").Value = 1999.5
$ws.Range("N51").Value = -2967.5

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H70").Value = 1899.7142
$ws.Range("I70").Value = 1421.75
$ws.Range("J70").Value = 2537
$ws.Range("K70").Value = 4265.25
$ws.Range("L70").Value = 7611
$ws.Range("M70").Value = -3995.25
$ws.Range("N70").Value = -8151

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H73").Value = 1899.7142
$ws.Range("I73").Value = 1421.75
$ws.Range("J73").Value = 2537
$ws.Range("K73").Value = 4265.25
$ws.Range("L73").Value = 7611
$ws.Range("M73").Value = -3329.25
$ws.Range("N73").Value = -9483

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H121").Value = 3831.3333
$ws.Range("J121").Value = 3831.3333
$ws.Range("L121").Value = 11493.9999
$ws.Range("N121").Value = -14987.9999

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H131").Value = 1252498.2
$ws.Range("I131").Value = 1667332.6
$ws.Range("K131").Value = 5001997.800000001
$ws.Range("M131").Value = -4996957.800000001

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H137").Value = 9498.629999999999
$ws.Range("I137").Value = 1839.1818
$ws.Range("J137").Value = 43200.2
$ws.Range("K137").Value = 5517.5454
$ws.Range("L137").Value = 129600.6
$ws.Range("M137").Value = -2967.5454
$ws.Range("N137").Value = -134700.6

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 3102.75
$ws.Range("I45").Value = 3117.4285
$ws.Range("K45").Value = 3117.4285
$ws.Range("M45").Value = -2740.4285

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 1868682.4
$ws.Range("I61").Value = 2899.2856
$ws.Range("J61").Value = 3174730.5
$ws.Range("K61").Value = 2899.2856
$ws.Range("L61").Value = 3174730.5
$ws.Range("M61").Value = -2687.2856
$ws.Range("N61").Value = -3175154.5

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H74").Value = 22542.607
$ws.Range("I74").Value = 1654.9048
$ws.Range("K74").Value = 1654.9048
$ws.Range("M74").Value = -780.9048

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H77").Value = 22542.607
$ws.Range("I77").Value = 1654.9048
$ws.Range("K77").Value = 8274.523999999999
$ws.Range("M77").Value = -3906.523999999999

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H132").Value = 6028969
$ws.Range("I132").Value = 4325.5
$ws.Range("K132").Value = 12976.5
$ws.Range("M132").Value = -10446.5

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H136").Value = 1868682.4
$ws.Range("I136").Value = 2899.2856
$ws.Range("J136").Value = 3174730.5
$ws.Range("K136").Value = 8697.856800000001
$ws.Range("L136").Value = 9524191.5
$ws.Range("M136").Value = -6147.856800000001
$ws.Range("N136").Value = -9529291.5

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 10443172
$ws.Range("I20").Value = 25649378
$ws.Range("J20").Value = 38926
$ws.Range("K20").Value = 25649378
$ws.Range("L20").Value = 38926
$ws.Range("M20").Value = -25649131
$ws.Range("N20").Value = -39420

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 17666.053
$ws.Range("J31").Value = 40607.5
$ws.Range("L31").Value = 40607.5
$ws.Range("N31").Value = -41197.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H34").Value = 17666.053
$ws.Range("J34").Value = 40607.5
$ws.Range("L34").Value = 40607.5
$ws.Range("N34").Value = -41011.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H93").Value = 11666
$ws.Range("I93").Value = 11666
$ws.Range("K93").Value = 11666
$ws.Range("M93").Value = -9794

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H99").Value = 7564175
$ws.Range("I99").Value = 5089964
$ws.Range("J99").Value = 13337333
$ws.Range("K99").Value = 5089964
$ws.Range("L99").Value = 13337333
$ws.Range("M99").Value = -5088466
$ws.Range("N99").Value = -13340329

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H126").Value = 7564175
$ws.Range("I126").Value = 5089964
$ws.Range("J126").Value = 13337333
$ws.Range("K126").Value = 15269892
$ws.Range("L126").Value = 40011999
$ws.Range("M126").Value = -15267422
$ws.Range("N126").Value = -40016939

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H8").Value = 62500040
$ws.Range("I8").Value = 62500040
$ws.Range("K8").Value = 187500120
$ws.Range("M8").Value = -187499981

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H31").Value = 6999.5
$ws.Range("I31").Value = 6999.5
$ws.Range("J31").Value = 0
$ws.Range("K31").Value = 20998.5
$ws.Range("L31").Value = 0
$ws.Range("M31").Value = -20710.5
$ws.Range("N31").ClearContents()

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H140").Value = 1638.5294
$ws.Range("I140").Value = 1638.5294
$ws.Range("J140").Value = 0
$ws.Range("K140").Value = 4915.5882
$ws.Range("L140").Value = 0
$ws.Range("M140").Value = 264.4117999999999
$ws.Range("N140").ClearContents()

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H141").Value = 3734.4614
$ws.Range("I141").Value = 2595.2727
$ws.Range("K141").Value = 7785.8181
$ws.Range("M141").Value = -2605.8181

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H5").Value = 2000
$ws.Range("J5").Value = 500
$ws.Range("L5").Value = 500
$ws.Range("N5").Value = -724

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H64").Value = 0
$ws.Range("J64").Value = 0
$ws.Range("L64").Value = 0
$ws.Range("N64").ClearContents()

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H67").Value = 0
$ws.Range("J67").Value = 0
$ws.Range("L67").Value = 0
$ws.Range("N67").ClearContents()

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 7140.0835
$ws.Range("I70").Value = 4767.5713
$ws.Range("J70").Value = 10461.6
$ws.Range("K70").Value = 4767.5713
$ws.Range("L70").Value = 10461.6
$ws.Range("M70").Value = -4497.5713
$ws.Range("N70").Value = -11001.6

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H73").Value = 7140.0835
$ws.Range("I73").Value = 4767.5713
$ws.Range("J73").Value = 10461.6
$ws.Range("K73").Value = 4767.5713
$ws.Range("L73").Value = 10461.6
$ws.Range("M73").Value = -3831.5713
$ws.Range("N73").Value = -12333.6

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 685143.25
$ws.Range("I132").Value = 3684.389
$ws.Range("J132").Value = 2437466
$ws.Range("K132").Value = 11053.167
$ws.Range("L132").Value = 7312398
$ws.Range("M132").Value = -8523.167000000001
$ws.Range("N132").Value = -7317458

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H4").Value = 13504.5
$ws.Range("I4").Value = 13504.5
$ws.Range("J4").Value = 0
$ws.Range("K4").Value = 13504.5
$ws.Range("L4").Value = 0
$ws.Range("M4").Value = -13391.5
$ws.Range("N4").ClearContents()

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 2076898.6
$ws.Range("I7").Value = 2807274.8
$ws.Range("K7").Value = 2807274.8
$ws.Range("M7").Value = -2807162.8

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H28").Value = 13504.5
$ws.Range("I28").Value = 13504.5
$ws.Range("J28").Value = 0
$ws.Range("K28").Value = 13504.5
$ws.Range("L28").Value = 0
$ws.Range("M28").Value = -13272.5
$ws.Range("N28").ClearContents()

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H37").Value = 13504.5
$ws.Range("I37").Value = 13504.5
$ws.Range("J37").Value = 0
$ws.Range("K37").Value = 13504.5
$ws.Range("L37").Value = 0
$ws.Range("M37").Value = -13397.5
$ws.Range("N37").ClearContents()

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H68").Value = 27500
$ws.Range("J68").Value = 5000
$ws.Range("L68").Value = 5000
$ws.Range("N68").Value = -6498

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H71").Value = 27500
$ws.Range("J71").Value = 5000
$ws.Range("L71").Value = 25000
$ws.Range("N71").Value = -32488

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H126").Value = 2076898.6
$ws.Range("I126").Value = 2807274.8
$ws.Range("K126").Value = 8421824.399999999
$ws.Range("M126").Value = -8419354.399999999

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H136").Value = 1434983
$ws.Range("J136").Value = 2492988.5
$ws.Range("L136").Value = 7478965.5
$ws.Range("N136").Value = -7484065.5

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H25").Value = 0
$ws.Range("J25").Value = 0
$ws.Range("L25").Value = 0
$ws.Range("N25").ClearContents()

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H29").Value = 8548.5
$ws.Range("J29").Value = 8548.5
$ws.Range("L29").Value = 8548.5
$ws.Range("N29").Value = -9128.5

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H104").Value = 30700
$ws.Range("J104").Value = 30700
$ws.Range("L104").Value = 30700
$ws.Range("N104").Value = -37688

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H113").Value = 2352.7297
$ws.Range("J113").Value = 2535.8333
$ws.Range("L113").Value = 7607.499899999999
$ws.Range("N113").Value = -11947.4999

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 356505.22
$ws.Range("J132").Value = 915483.3
$ws.Range("L132").Value = 2746449.9
$ws.Range("N132").Value = -2751509.9
